# Apply the "new exhibition row inserted at row 7" update to both the
# "展览" (index 1) and "全部类型" (index 4) worksheets. Sheet "演出" (2) and
# "本地生活" (3) are untouched (they only contain the header row).

$wb = $excel.ActiveWorkbook

$sheetIndexes = @(1, 4)

foreach ($si in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($si)

    # 1) Overall "want to go" total (row 2, column F) ticks up by one.
    $ws.Range("F2").Value2 = 1165

    # 2) Insert a new row at position 7; this shifts the previous rows
    #    7..22 down to 8..23 (formatting/values carried with them).
    $ws.Rows.Item(7).Insert()

    # The freshly inserted row 7 only got a border-less clone style on A7;
    # repair it by pulling the real "index column" format from the row
    # that was just pushed down to row 8 (bold/center/bordered, style s=1).
    $ws.Range("A8").Copy()
    $ws.Range("A7").PasteSpecial(-4122)   # xlPasteFormats
    $ws.Application.CutCopyMode = 0

    # 3) Populate the new row 7 with the new event's data.
    $ws.Range("A7").Value2 = 6

    # B7 ("2024-03-31") would otherwise be auto-coerced into a real Excel
    # date serial. Write it with a forcing leading apostrophe (quote
    # prefix) so it lands as plain text, then strip the resulting
    # quote-prefix style by re-pasting the (unstyled) format from the
    # neighbouring text cell C7 — leaving a clean, un-styled text value
    # exactly like the rest of column B.
    $ws.Range("B7").Value2 = "'2024-03-31"
    $ws.Range("C7").Value2 = "张家港·META萌圆饿了"
    $ws.Range("C7").Copy()
    $ws.Range("B7").PasteSpecial(-4122)   # xlPasteFormats
    $ws.Application.CutCopyMode = 0

    $ws.Range("D7").Value2 = "大新镇人民路18号 新香苑宴会厅"
    $ws.Range("E7").Value2 = "2024.03.31 10:00-03.31 17:00"
    $ws.Range("F7").Value2 = 0
    $ws.Range("G7").Value2 = 30
    $ws.Range("H7").Value2 = "https://show.bilibili.com/platform/detail.html?id=82407"
    $ws.Range("I7").Value2 = "//i0.hdslb.com/bfs/openplatform/202403/pxZkSPcL1709707210340.jpeg"

    # 4) A handful of the shifted rows carry slightly refreshed "want to
    #    go" counts (F column) alongside the pure positional shift.
    $ws.Range("F8").Value2 = 11335    # was 11330 (Anime LIVE)
    $ws.Range("F10").Value2 = 88      # was 89 (寒假动漫展宅舞比赛-CF01)
    $ws.Range("F15").Value2 = 12320   # was 12316 (理想乡动漫游戏展)
    $ws.Range("F16").Value2 = 12973   # was 12969 (I COME ACG 动漫品牌博览会)
    $ws.Range("F23").Value2 = 89      # was 88 (Come in joy动漫国潮文化节)
}

Write-Output "done"
